$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

# 1. Remove the "PayPoint" asset row (HR_FTAE_Aurion_PayPoint_FilePath / PayPointFile)
$ws.Rows.Item(5).Delete()

# 2. Repurpose the old "Dev and Test email" asset (now shifted up to row 9) into the
#    new "DefaultTaskEmailAddress" asset. Set column B (key) before column A (name)
#    so the shared-string table is built in the same order as the authored workbook.
$ws.Cells.Item(9, 2).Value2 = "HR_FTAE_DefaultTaskEmailAddress    "
$ws.Cells.Item(9, 1).Value2 = "DefaultTaskEmailAddress    "
$ws.Cells.Item(9, 3).Value2 = "Email address for task if no PST match found"

# 3. Re-sort the Assets table (rows below the header) alphabetically by column B,
#    the same way Excel's Data > Sort records it (range padded out to row 1001).
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B1001")) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:C1001"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# 4. Restore the selection left behind in the saved file.
$ws.Range("C6").Select() | Out-Null
